$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.01297935103244838
$ws.Range("C2").Value = 0.01238856843209229
$ws.Range("D2").Value = 0.01272131147540984
$ws.Range("E2").Value = 0.01187196641742096
$ws.Range("F2").Value = 0.01301419135439147
$ws.Range("G2").Value = 0.01175003282132073
$ws.Range("H2").Value = 0.01224222585924714
$ws.Range("I2").Value = 0.01205766710353866
$ws.Range("J2").Value = 0.01290364839195651
$ws.Range("K2").Value = 0.01284319507240679
$ws.Range("L2").Value = 0.01166601127277494
$ws.Range("M2").Value = 0.01160807974816369
$ws.Range("N2").Value = 0.01251638269986894
$ws.Range("O2").Value = 0.01155538047403322
$ws.Range("P2").Value = 0.01160199265862611
$ws.Range("Q2").Value = 0.01210970740328599
$ws.Range("R2").Value = 0.01137633984349313
$ws.Range("S2").Value = 0.01192191798768505
$ws.Range("T2").Value = 0.01218872870249017
$ws.Range("U2").Value = 0.01322595429843515
$ws.Range("V2").Value = 0.01205766710353866
$ws.Range("W2").Value = 0.01167825744652933
$ws.Range("X2").Value = 0.0111599816188538
$ws.Range("Y2").Value = 0.01226792626123466
$ws.Range("Z2").Value = 0.012
$ws.Range("AA2").Value = 0.01187897880160137
$ws.Range("AB2").Value = 0.01175311884438608
$ws.Range("AC2").Value = 0.01262015301118159
$ws.Range("AD2").Value = 0.01251392255781956
$ws.Range("AE2").Value = 0.01290026848274507
$ws.Range("AF2").Value = 0.01253280839895013
$ws.Range("AG2").Value = 0.01296660117878192
$ws.Range("AH2").Value = 0.01178010471204189
$ws.Range("AI2").Value = 0.01261850277868585
$ws.Range("AJ2").Value = 0.0133604034317899
$ws.Range("AK2").Value = 0.01230688661953391
$ws.Range("AL2").Value = 0.01277934333835769
$ws.Range("AM2").Value = 0.01206161914126516
$ws.Range("AN2").Value = 0.01198271346254584
$ws.Range("AO2").Value = 0.01154401154401154
$ws.Range("AP2").Value = 0.0127826941986234
$ws.Range("AQ2").Value = 0.01367264163286667
$ws.Range("AR2").Value = 0.01257532093266964
$ws.Range("AS2").Value = 0.01256462273411426
$ws.Range("AT2").Value = 0.01244514311914587
$ws.Range("AU2").Value = 0.01167825744652933
$ws.Range("AV2").Value = 0.01258438749426493
$ws.Range("AW2").Value = 0.01359654856844032
$ws.Range("AX2").Value = 0.01349138777916039
$ws.Range("AY2").Value = 0.0130167451596023
$ws.Range("AZ2").Value = 0.0123542734169672
$ws.Range("BA2").Value = 0.000603791569630046
# Row 3
$ws.Range("B3").Value = 0.3245901639344262
$ws.Range("C3").Value = 0.3103448275862069
$ws.Range("D3").Value = 0.3164763458401305
$ws.Range("E3").Value = 0.2986798679867987
$ws.Range("F3").Value = 0.3241042345276873
$ws.Range("G3").Value = 0.2905844155844156
$ws.Range("H3").Value = 0.3164128595600677
$ws.Range("I3").Value = 0.3041322314049587
$ws.Range("J3").Value = 0.3172302737520129
$ws.Range("K3").Value = 0.3197389885807504
$ws.Range("L3").Value = 0.3016949152542373
$ws.Range("M3").Value = 0.2974789915966387
$ws.Range("N3").Value = 0.308562197092084
$ws.Range("O3").Value = 0.2880523731587561
$ws.Range("P3").Value = 0.2915980230642504
$ws.Range("Q3").Value = 0.3027823240589198
$ws.Range("R3").Value = 0.2957264957264957
$ws.Range("S3").Value = 0.2940226171243942
$ws.Range("T3").Value = 0.3136593591905565
$ws.Range("U3").Value = 0.3300653594771242
$ws.Range("V3").Value = 0.3108108108108108
$ws.Range("W3").Value = 0.2986577181208054
$ws.Range("X3").Value = 0.2866779089376054
$ws.Range("Y3").Value = 0.3101160862354892
$ws.Range("Z3").Value = 0.3029801324503311
$ws.Range("AA3").Value = 0.3052276559865092
$ws.Range("AB3").Value = 0.3028764805414552
$ws.Range("AC3").Value = 0.3254637436762226
$ws.Range("AD3").Value = 0.310064935064935
$ws.Range("AE3").Value = 0.3234811165845649
$ws.Range("AF3").Value = 0.3141447368421053
$ws.Range("AG3").Value = 0.3219512195121951
$ws.Range("AH3").Value = 0.3010033444816054
$ws.Range("AI3").Value = 0.3356521739130435
$ws.Range("AJ3").Value = 0.3417085427135678
$ws.Range("AK3").Value = 0.3128119800332779
$ws.Range("AL3").Value = 0.3207236842105263
$ws.Range("AM3").Value = 0.300163132137031
$ws.Range("AN3").Value = 0.2990196078431372
$ws.Range("AO3").Value = 0.284329563812601
$ws.Range("AP3").Value = 0.3170731707317073
$ws.Range("AQ3").Value = 0.3370967741935484
$ws.Range("AR3").Value = 0.3127035830618892
$ws.Range("AS3").Value = 0.320534223706177
$ws.Range("AT3").Value = 0.319327731092437
$ws.Range("AU3").Value = 0.2951907131011609
$ws.Range("AV3").Value = 0.318407960199005
$ws.Range("AW3").Value = 0.3280757097791798
$ws.Range("AX3").Value = 0.3433333333333333
$ws.Range("AY3").Value = 0.32569558101473
$ws.Range("AZ3").Value = 0.3114248043724379
$ws.Range("BA3").Value = 0.01424057628614986
# Row 4
$ws.Range("B4").Value = 0.02496060510557832
$ws.Range("C4").Value = 0.02382603214623385
$ws.Range("D4").Value = 0.02445943390279266
$ws.Range("E4").Value = 0.0228362351753722
$ws.Range("F4").Value = 0.02502357749135492
$ws.Range("G4").Value = 0.02258675078864353
$ws.Range("H4").Value = 0.02357241900920207
$ws.Range("I4").Value = 0.02319571383548692
$ws.Range("J4").Value = 0.02479859013091642
$ws.Range("K4").Value = 0.02469446894292554
$ws.Range("L4").Value = 0.02246340232205957
$ws.Range("M4").Value = 0.0223442529823897
$ws.Range("N4").Value = 0.02405693053718748
$ws.Range("O4").Value = 0.02221941674031057
$ws.Range("P4").Value = 0.0223160814473933
$ws.Range("Q4").Value = 0.02328801611278953
$ws.Range("R4").Value = 0.02190982776089159
$ws.Range("S4").Value = 0.02291469940195153
$ws.Range("T4").Value = 0.02346559010912761
$ws.Range("U4").Value = 0.02543279823733081
$ws.Range("V4").Value = 0.0232147363108756
$ws.Range("W4").Value = 0.02247758555373153
$ws.Range("X4").Value = 0.02148363452546443
$ws.Range("Y4").Value = 0.02360217089486305
$ws.Range("Z4").Value = 0.02308565661662672
$ws.Range("AA4").Value = 0.02286797220467467
$ws.Range("AB4").Value = 0.02262815245559699
$ws.Range("AC4").Value = 0.02429812413445802
$ws.Range("AD4").Value = 0.02405693053718748
$ws.Range("AE4").Value = 0.02481108312342569
$ws.Range("AF4").Value = 0.02410398788490661
$ws.Range("AG4").Value = 0.02492917847025496
$ws.Range("AH4").Value = 0.02267288071545535
$ws.Range("AI4").Value = 0.02432262129804663
$ws.Range("AJ4").Value = 0.0257153661918568
$ws.Range("AK4").Value = 0.0236820558039932
$ws.Range("AL4").Value = 0.02457931556059747
$ws.Range("AM4").Value = 0.0231913284597933
$ws.Range("AN4").Value = 0.02304205489801058
$ws.Range("AO4").Value = 0.02218720453829184
$ws.Range("AP4").Value = 0.02457466918714556
$ws.Range("AQ4").Value = 0.02627939142461964
$ws.Range("AR4").Value = 0.0241783150736683
$ws.Range("AS4").Value = 0.02418136020151133
$ws.Range("AT4").Value = 0.02395662589837348
$ws.Range("AU4").Value = 0.02246765541180183
$ws.Range("AV4").Value = 0.02421185372005044
$ws.Range("AW4").Value = 0.02611097162942505
$ws.Range("AX4").Value = 0.02596256852983805
$ws.Range("AY4").Value = 0.0250330209447135
$ws.Range("AZ4").Value = 0.02376546628758393
$ws.Range("BA4").Value = 0.001155664846199266
# Row 5
$ws.Range("B5").Value = 0.04679860302677532
$ws.Range("C5").Value = 0.04679860302677532
$ws.Range("D5").Value = 0.04679860302677532
$ws.Range("E5").Value = 0.04633294528521537
$ws.Range("F5").Value = 0.04633294528521537
$ws.Range("G5").Value = 0.04633294528521537
$ws.Range("H5").Value = 0.04679860302677532
$ws.Range("I5").Value = 0.04679860302677532
$ws.Range("J5").Value = 0.04679860302677532
$ws.Range("K5").Value = 0.0470314318975553
$ws.Range("L5").Value = 0.04726426076833527
$ws.Range("N5").Value = 0.04633294528521537
$ws.Range("O5").Value = 0.0470314318975553
$ws.Range("P5").Value = 0.04679860302677532
$ws.Range("R5").Value = 0.04679860302677532
$ws.Range("T5").Value = 0.0470314318975553
$ws.Range("W5").Value = 0.0470314318975553
$ws.Range("X5").Value = 0.04726426076833527
$ws.Range("AA5").Value = 0.0470314318975553
$ws.Range("AB5").Value = 0.04679860302677532
$ws.Range("AC5").Value = 0.04679860302677532
$ws.Range("AE5").Value = 0.04656577415599535
$ws.Range("AF5").Value = 0.04656577415599535
$ws.Range("AG5").Value = 0.04679860302677532
$ws.Range("AH5").Value = 0.04679860302677532
$ws.Range("AJ5").Value = 0.0470314318975553
$ws.Range("AL5").Value = 0.0470314318975553
$ws.Range("AM5").Value = 0.04656577415599535
$ws.Range("AO5").Value = 0.04633294528521537
$ws.Range("AP5").Value = 0.04633294528521537
$ws.Range("AS5").Value = 0.04679860302677532
$ws.Range("AT5").Value = 0.04633294528521537
$ws.Range("AU5").Value = 0.0470314318975553
$ws.Range("AV5").Value = 0.0470314318975553
$ws.Range("AZ5").Value = 0.04680325960419091
$ws.Range("BA5").Value = 0.0002654249126891741
